$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41 (pushes existing rows 41-61 down to 42-62)
$ws.Rows.Item(41).Insert()

# Populate the new row 41 with the new weekly data point
$ws.Range("A41").Value = 1
$ws.Range("B41").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C41").Value = "Arica y Parinacota"
$ws.Range("D41").Value = 44726
$ws.Range("E41").Value = 15
$ws.Range("F41").Value = 100112031
$ws.Range("G41").Value = "Poroto verde"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 1500
$ws.Range("K41").Value = 800
$ws.Range("L41").Value = 900
$ws.Range("M41").Value = 850
$ws.Range("N41").Value = "$/kilo"
$ws.Range("O41").Value = "Región de Arica y Parinacota"
$ws.Range("P41").Value = 850
$ws.Range("Q41").Value = 1
$ws.Range("R41").Value = "Hortaliza"
